$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 228, shifting existing rows 228-288 down to 229-289.
$ws.Rows.Item(228).Insert()

# Populate the new row 228 with the same static columns as the surrounding
# records, plus the new date (2022-03-22) and volume (160).
$ws.Cells.Item(228, 1).Value = 3
$ws.Cells.Item(228, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(228, 3).Value = "Coquimbo"
$ws.Cells.Item(228, 4).Value = Get-Date -Year 2022 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(228, 5).Value = 5
$ws.Cells.Item(228, 6).Value = 100112039
$ws.Cells.Item(228, 7).Value = "Ciboulette"
$ws.Cells.Item(228, 8).Value = "Sin especificar"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 160
$ws.Cells.Item(228, 11).Value = 1500
$ws.Cells.Item(228, 12).Value = 1500
$ws.Cells.Item(228, 13).Value = 1500
$ws.Cells.Item(228, 14).Value = "`$/docena de atados"
$ws.Cells.Item(228, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(228, 16).Value = 500
$ws.Cells.Item(228, 17).Value = 3
$ws.Cells.Item(228, 18).Value = "Hortaliza"
